$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.537.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.290.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.35%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.564"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.81%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.643.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.289.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.805"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.530.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0932"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0789"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.40%  "

$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +20.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0301"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.818.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "86.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +19.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.193"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "72.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.515.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.18%  "
